$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 6 -------------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A6").Value = "Ik heb nog geen geld terug."
$ws.Range("B6").Value = "mailmind.test@zohomail.eu"
$ws.Range("C6").Value = "Testmail #1: Ik heb nog geen geld terug."
$ws.Range("D6").Value = "Retour / Terugbetaling"
$ws.Range("E6").Value = "Beste klant,`nBedankt voor uw bericht. Om uw vraag te kunnen beantwoorden, heb ik wat meer informatie nodig. Kunt u mij de volgende gegevens bezorgen:`n1. Het ordernummer van uw aankoop.`n2. De datum waarop de terugbetaling is aangevraagd.`n3. Eventuele referentienummers of andere relevante details.`nMet deze gegevens kan ik uw zaak verder onderzoeken en u van passende ondersteuning voorzien.`nMet vriendelijke groet,`n[Naam] `nKlantenservice Team"
$ws.Range("F6").Value = "2025-08-04 20:07:30"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"

# --- Expand conditional-formatting ranges from row 5 to row 6 ----------
$cols = "D", "G", "H", "I", "J"
foreach ($col in $cols) {
    $fcs = $ws.Range($col + "2").FormatConditions
    $newRange = $ws.Range($col + "2:" + $col + "6")
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump "Retour / Terugbetaling" count to 2 ---------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 2
